$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.049.52"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "1.922.37"
$ws.Range("E3").Value = "  +1.97%  "

$ws.Range("E4").Value = "  +0.55%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.60%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4604"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3821"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07765"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9782"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.87%  "

$ws.Range("D12").Value = "1.914.93"
$ws.Range("E12").Value = "  +1.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.685"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.959"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07069"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009515"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.04%  "

$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").Value = "29.110.34"
$ws.Range("E21").Value = "  +1.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.349"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.078"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.662"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "118.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.840"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09338"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8570"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.121"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.244"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.025"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.55%  "

$ws.Range("B35").Value = "TrustWalletToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.162"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.00%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05683"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.67%  "

$ws.Range("E37").Value = "  +0.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02049"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.141"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +17.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.498"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5520"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1755"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.309"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000002829"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.200"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5210"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.31%  "

$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06931"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.766"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("E51").Value = "  +0.61%  "
